$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# New handoff identifiers / timestamps (this run's "Generate Report for Handoff")
$oldGuid = "412ecf0e-d941-4ca7-82c7-1e0e9d269306"
$newGuid = "174a2b1c-331b-42b8-8515-e6d9ffe67458"
$oldHash = "9fb6075c468c3adea63edbbf88d4e2ab909970e1"
$newHash = "46b78910d0359d63da7fa33daccc84859dc69c5e"

$newMdName = "$newGuid.md"
$newZhXlf  = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf  = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet: File Name (A2) + Latest Handoff Date (D2)
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = $newMdName
$ws1.Range("D2").Value = "2016-03-22 19:02:35"

$link = $ws1.Range("A2")
$link.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($link, "https://github.com/OpenLocalizationTest/oltest/blob/b3630b86daafb6d62561b3c92dc111a5c295b873/e2e/$oldGuid.md", "", "", $newMdName)

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name (A2), Latest Handoff File (D2), Latest Handoff Datetime (E2)
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = $newMdName
$ws2.Range("D2").Value = $newZhXlf
$ws2.Range("E2").Value = "2016-03-22 19:02:29"

$linkA = $ws2.Range("A2")
$linkD = $ws2.Range("D2")
$linkA.Hyperlinks.Delete()
$linkD.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($linkA, "https://github.com/OpenLocalizationTest/oltest/blob/b3630b86daafb6d62561b3c92dc111a5c295b873/e2e/$oldGuid.md", "", "", $newMdName)
$ws2.Hyperlinks.Add($linkD, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fafe73fe6429da5d177debba24803ce1e1bf3e24/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf", "", "", $newZhXlf)

# ---------------------------------------------------------------------------
# de-de sheet: Source File Name (A2), Latest Handoff File (D2), Latest Handoff Datetime (E2)
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = $newMdName
$ws3.Range("D2").Value = $newDeXlf
$ws3.Range("E2").Value = "2016-03-22 19:02:35"

$linkA3 = $ws3.Range("A2")
$linkD3 = $ws3.Range("D2")
$linkA3.Hyperlinks.Delete()
$linkD3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($linkA3, "https://github.com/OpenLocalizationTest/oltest/blob/b3630b86daafb6d62561b3c92dc111a5c295b873/e2e/$oldGuid.md", "", "", $newMdName)
$ws3.Hyperlinks.Add($linkD3, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7122abd85d041b72824065455f4849e7940002d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf", "", "", $newDeXlf)
